$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "last updated" timestamp label (22:22 -> 22:52)
$ws.Range("A1").Value = "Datos actualizados a 19 de Abril de 2020 a las 22:52"

# 2/3/4. Brasil overtakes Belgica: Brasil moves up to row 14 with refreshed
# numbers, Belgica drops to row 15 keeping its previous numbers.
$ws.Range("A14").Value = "Brasil"
$ws.Range("B14").Value = 38654
$ws.Range("C14").Value = 1932
$ws.Range("D14").Value = 14026
$ws.Range("E14").Value = 22166
$ws.Range("F14").Value = 6634
$ws.Range("G14").Value = 101
$ws.Range("H14").Value = 2462

$ws.Range("A15").Value = "Belgica"
$ws.Range("B15").Value = 38496
$ws.Range("C15").Value = 1313
$ws.Range("D15").Value = 8757
$ws.Range("E15").Value = 24056
$ws.Range("F15").Value = 1081
$ws.Range("G15").Value = 230
$ws.Range("H15").Value = 5683

# 5/6/7. San Martin (Parte Holandesa) overtakes Sudan: San Martin moves up to
# row 150 with refreshed numbers, Sudan drops to row 151 keeping its numbers.
$ws.Range("A150").Value = "San Martin (Parte Holandesa)"
$ws.Range("B150").Value = 67
$ws.Range("C150").Value = 3
$ws.Range("D150").Value = 12
$ws.Range("E150").Value = 45
$ws.Range("F150").Value = 6
$ws.Range("G150").Value = 1
$ws.Range("H150").Value = 10

$ws.Range("A151").Value = "Sudan"
$ws.Range("B151").Value = 66
$ws.Range("C151").Value = 0
$ws.Range("D151").Value = 6
$ws.Range("E151").Value = 50
$ws.Range("F151").Value = 0
$ws.Range("G151").Value = 0
$ws.Range("H151").Value = 10

# 8. Refresh Estados Unidos row totals
$ws.Range("B4").Value = 761379
$ws.Range("C4").Value = 22587
$ws.Range("D4").Value = 69929
$ws.Range("E4").Value = 651031
$ws.Range("F4").Value = 13556
$ws.Range("G4").Value = 1405
$ws.Range("H4").Value = 40419

# 9. Republica Dominicana: only "Muertes hoy" changes
$ws.Range("F48").Value = 114

# 10. Tunez row refresh
$ws.Range("B85").Value = 879
$ws.Range("C85").Value = 13
$ws.Range("D85").Value = 43
$ws.Range("E85").Value = 798
$ws.Range("F85").Value = 33
$ws.Range("G85").Value = 1
$ws.Range("H85").Value = 38

# 11. Niger row refresh
$ws.Range("B95").Value = 648
$ws.Range("C95").Value = 9
$ws.Range("D95").Value = 117
$ws.Range("E95").Value = 511
$ws.Range("F95").Value = 0
$ws.Range("G95").Value = 1
$ws.Range("H95").Value = 20

# 12. Estado de Palestina row refresh
$ws.Range("B105").Value = 437
$ws.Range("C105").Value = 19
$ws.Range("D105").Value = 71
$ws.Range("E105").Value = 363
$ws.Range("F105").Value = 0
$ws.Range("G105").Value = 1
$ws.Range("H105").Value = 3
